$wb = $excel.ActiveWorkbook

foreach ($ws in @($wb.Worksheets.Item(1), $wb.Worksheets.Item(2))) {
    # Rename the surviving "adult5" columns to the merged "adult" label.
    # Set adultM (col I) before adultF (col E) so the shared-string table
    # picks up the same ordering as the target workbook.
    $ws.Range("I1").Value2 = "adultM"
    $ws.Range("E1").Value2 = "adultF"

    # Drop the "adult3" columns entirely (D = adult3F, H = adult3M after
    # the prior deletion shifts columns left).
    $ws.Columns("D:D").Delete()
    $ws.Columns("G:G").Delete()

    # Re-apply the row total as one shared formula across H2:H13 (column
    # deletion above turns it back into per-cell formulas), matching the
    # original "sum across the row" layout one column to the left.
    $ws.Range("H2:H13").Formula = "=SUM(B2:G2)"

    $ws.Range("A1:H1").Select()
}
